$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the Sep 9 2023 GitHub Actions data refresh.
# Some Price values are plain decimal-looking numbers ("216.11", "63.00", ...). Excel would
# auto-convert those to numeric cells (dropping the trailing zero / changing type), so for
# those specific cells we force the cell to Text format before writing the value, matching
# the original inline-string text cells.

$ws.Range("D2").Value = "25.982.77"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "1.641.42"
$ws.Range("E3").Value = "  +0.94%  "

$ws.Range("E4").Value = "  +0.52%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.11"

$ws.Range("E6").Value = "  +1.39%  "

$ws.Range("E8").Value = "  +0.64%  "

$ws.Range("E9").Value = "  +1.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.65"
$ws.Range("E10").Value = "  +0.11%  "

$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("E12").Value = "  +1.09%  "

$ws.Range("D13").Value = "1.867.98"
$ws.Range("E13").Value = "  +0.94%  "

$ws.Range("D14").Value = "1.642.34"
$ws.Range("E14").Value = "  +0.78%  "

$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("D16").Value = "0.0₃0766"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.00"
$ws.Range("E17").Value = "  +0.71%  "

$ws.Range("D18").Value = "25.955.33"
$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.01"
$ws.Range("E19").Value = "  +0.50%  "

$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.93"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.24"
$ws.Range("E23").Value = "  +0.32%  "

$ws.Range("E24").Value = "  +6.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.82"
$ws.Range("E25").Value = "  +1.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.46"
$ws.Range("E26").Value = "  +1.46%  "

$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.92"
$ws.Range("E28").Value = "  +1.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.53"
$ws.Range("E29").Value = "  +0.60%  "

$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0500"
$ws.Range("E31").Value = "  +0.52%  "

$ws.Range("E32").Value = "  -0.83%  "

$ws.Range("E33").Value = "  +1.20%  "

$ws.Range("E34").Value = "  -2.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").Value = "  +2.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.905"
$ws.Range("E36").Value = "  +0.33%  "

$ws.Range("D37").Value = "1.132.81"
$ws.Range("E37").Value = "  +0.31%  "

$ws.Range("E38").Value = "  -1.04%  "

$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("E40").Value = "  +0.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.53"
$ws.Range("E41").Value = "  +1.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.27"
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "1.776.61"
$ws.Range("E44").Value = "  +0.84%  "

$ws.Range("E45").Value = "  +3.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.56"
$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("E48").Value = "  +0.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.73"
$ws.Range("E49").Value = "  +2.07%  "

$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("E51").Value = "  +0.18%  "
